$d = $word.ActiveDocument

# Change 1: extend first SmartNodes paragraph with Russian payout sentence
$d.Content.Find.Execute(
    "SmartNodes are paid at a rate of 10 nodes every other block, which means a payout will be 2% of the Block Reward for every Node.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SmartNodes получают выплату следующим образом: 1 нода каждый второй блок, что означает, что выплата будет равна 2% от награды за блок для каждой ноды.",
    2)

# Change 2: InstantPay bullet translated to Russian
$d.Content.Find.Execute(
    "InstantPay (Instant Transactions): Allows for SmartCash transactions to be locked in about a second. No risk of double spending a transaction, so the receiver can trust that transaction immediately.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "InstantPay (мгновенные транзакции): Позволяет практически мгновенно переводить средства любому пользователю. Нет риска двойной траты, поэтому получатель может сразу же доверять транзакции.",
    2)

# Change 3: SmartRewards bullet translated to Russian
$d.Content.Find.Execute(
    "SmartRewards: SmartRewards are calculated by the SmartNodes to allow for distribution to be handled automatically by the block rewards.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SmartRewards: SmartRewards вычисляется с помощью SmartNodes, которые позволяют осуществлять автоматическое распределение награды.",
    2)

# Change 4: "How much will I earn" paragraph translated to Russian
$d.Content.Find.Execute(
    "SmartNodes are paid at a rate of 1 node every other block, which means a payout will be 2% of the Block Reward for every Node",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SmartNodes получают выплату следующим образом: 1 нода каждый второй блок, что означает, что выплата будет равна 2% от награды за блок для каждой ноды.",
    2)
